{"js": "const results = context.document.body.search(\"surname\", { matchCase: false });\nresults.load(\"text,items\");\nawait context.sync();\nlet out = [];\nfor (const r of results.items) {\n  r.load(\"text\");\n}\nawait context.sync();\nfor (const r of results.items) {\n  out.push(r.text);\n}\nreturn JSON.stringify(out);\n", "ps1": "$d = $word.ActiveDocument\n$s = $d.Shapes(\"Text Box 10\")\n$tf = $s.TextFrame\n$tr = $tf.TextRange\nWrite-Output \"ParagraphsCount: $($tr.Paragraphs.Count)\"\nfor ($i=1; $i -le $tr.Paragraphs.Count; $i++) {\n    $p = $tr.Paragraphs.Item($i)\n    Write-Output \"$i => $($p.Range.Text)\"\n}\n"}
